$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the old CRC block (F1:K6, merged G1:J1) -------------------
$ws.Range("G1:J1").UnMerge()
$ws.Rows("1:6").Delete()

# --- Header row (row 8): dec / hex labels + merged blank area ---------
$ws.Range("L8").Value = "dec"
$ws.Range("M8").Value = "hex"
$ws.Range("M8:P8").Merge()

# --- Type / Address / Length rows (9-11) -------------------------------
$ws.Range("L9").Formula = "=HEX2DEC(P9)"
$ws.Range("P9").Value = 1
$ws.Range("Q9").Value = "Type"

$ws.Range("L10").Formula = "=HEX2DEC(P10)"
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = "Address"

$ws.Range("L11").Formula = "=HEX2DEC(P11)"
$ws.Range("P11").Value = 5
$ws.Range("Q11").Value = "Length"

# --- Payload checksum row (12) -----------------------------------------
$ws.Range("L12").Formula = "=SUM((Q15:Q24))"

# --- CRC result row (13) -------------------------------------------------
$ws.Range("L13").Formula = "=DEC2HEX(MOD(SUM(L9:L11,L12),256))"
$ws.Range("M13").Value = "CRC [hex]"

# --- Payload lookup table header (row 14, merged Q14:R14) --------------
$ws.Range("Q14").Value = "Payload "
$ws.Range("Q14:R14").Merge()

# --- Payload lookup table (rows 15-20) ----------------------------------
$ws.Range("Q15").Formula = "=HEX2DEC(R15)"
$ws.Range("R15").Value = "A1"

$ws.Range("Q16").Formula = "=HEX2DEC(R16)"
$ws.Range("R16").Value = "B2"

$ws.Range("Q17").Formula = "=HEX2DEC(R17)"
$ws.Range("R17").Value = "C3"

$ws.Range("Q18").Formula = "=HEX2DEC(R18)"
$ws.Range("R18").Value = "D4"

$ws.Range("Q19").Formula = "=HEX2DEC(R19)"
$ws.Range("R19").Value = "E5"

$ws.Range("Q20").Formula = "=HEX2DEC(R20)"
$ws.Range("R20").Value = "F6"

# --- Column widths for the new table ------------------------------------
$ws.Columns("P").ColumnWidth = 7.25
$ws.Columns("Q").ColumnWidth = 21.5

# --- View / selection ----------------------------------------------------
$ws.Range("O13").Select()
$excel.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 7

Write-Output "done"
